$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '41.840.28'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  -1.08%  '

# Row 3
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.480.54'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  -1.84%  '

# Row 4
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.00'
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  +0.06%  '

# Row 5
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '310.95'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  +0.17%  '

# Row 6
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '95.12'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  -2.91%  '

# Row 7
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.553'
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  -2.08%  '

# Row 8
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  +0.08%  '

# Row 9
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.510'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  -2.67%  '

# Row 10
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '34.17'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  -3.70%  '

# Row 11
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0788'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  -1.31%  '

# Row 12
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.108'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  +0.58%  '

# Row 13
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '7.04'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  -3.56%  '

# Row 14
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '2.865.87'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  -1.98%  '

# Row 15
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '2.513.78'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  +2.13%  '

# Row 16
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  -6.22%  '

# Row 17
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.791'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  -4.84%  '

# Row 18
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '41.865.80'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  -1.02%  '

# Row 19
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '6.39'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  -5.43%  '

# Row 20
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.0₃0923'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  -2.27%  '

# Row 21
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '11.73'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  -3.35%  '

# Row 22
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '68.93'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  +0.29%  '

# Row 23
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '237.44'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  -2.03%  '

# Row 24
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  -2.67%  '

# Row 25
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  -4.31%  '

# Row 26
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  +0.06%  '

# Row 27
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '24.72'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  -4.95%  '

# Row 28
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  -4.46%  '

# Row 29
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '9.78'
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  -3.11%  '

# Row 30
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '36.79'
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  -6.00%  '

# Row 31
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '154.93'
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  -1.60%  '

# Row 32
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '5.66'
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  +0.02%  '

# Row 33
$ws.Range('B33').NumberFormat = '@'
$ws.Range('B33').Value = 'WEMIXToken'
$ws.Range('C33').NumberFormat = '@'
$ws.Range('C33').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '2.61'
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  -0.38%  '

# Row 34
$ws.Range('B34').NumberFormat = '@'
$ws.Range('B34').Value = 'Hedera'
$ws.Range('C34').NumberFormat = '@'
$ws.Range('C34').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.0761'
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  -3.76%  '

# Row 35
$ws.Range('B35').NumberFormat = '@'
$ws.Range('B35').Value = 'ApeXProtocol'
$ws.Range('C35').NumberFormat = '@'
$ws.Range('C35').Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '2.49'
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  -11.63%  '

# Row 36
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  -4.43%  '

# Row 37
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '17.24'
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  -4.94%  '

# Row 38
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  -6.31%  '

# Row 39
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  -3.94%  '

# Row 40
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  -2.23%  '

# Row 41
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '4.06'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  -3.96%  '

# Row 42
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '21.39'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  -0.26%  '

# Row 43
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  -0.12%  '

# Row 44
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '2.003.32'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  +2.42%  '

# Row 45
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  -2.58%  '

# Row 46
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '3.10'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  -5.17%  '

# Row 47
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '8.73'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  -1.76%  '

# Row 48
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.723.56'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  -1.96%  '

# Row 49
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '77.51'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  -3.64%  '

# Row 50
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '70.12'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  -2.50%  '

# Row 51
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  -3.93%  '
